$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'298.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.14%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'31.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.14%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.147"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.63%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.07962"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'8.16%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'2.600"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'57.74%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'7.836"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.86%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'3.828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.27%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.9073"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.07%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.1731"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.50%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07239"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.12%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.08012"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.99%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.03022"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.47%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.09971"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.71%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.001496"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.68%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.005907"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.88%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'3.507"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.255"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.31%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'0.34%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.1328"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.42%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'4.609"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.19%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.1599"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.29%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04590"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.60%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.001259"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.61%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.004457"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.79%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.0001179"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.0003428"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'83.24%"
$ws.Range("E27").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.01839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'9.50%"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.04524"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.34%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.007021"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.90%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.1344"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.29%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.002238"
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "'-5.24%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00006472"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.44%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = 'BOLO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D47").Value = "'0.8206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'15.29%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = "'0.006196"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-39.33%"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = "Normal"
